$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F4, F5, F7, F9
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 469
$ws1.Range("F5").Value = 491
$ws1.Range("F7").Value = 2547
$ws1.Range("F9").Value = 6873

# Sheet "全部类型" (sheet4): update F4, F5, F9, F11
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 469
$ws4.Range("F5").Value = 491
$ws4.Range("F9").Value = 2547
$ws4.Range("F11").Value = 6873
